$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the refreshed crypto data.
# Leading apostrophe forces Excel to treat numeric-looking price strings as text,
# matching the original inline-string cell type; the cell Style is restored
# immediately afterward so no stray numeric formatting is left behind.

$dCell = $ws.Cells.Item(2, 4)
$dStyle = $dCell.Style
$dCell.Value = "'22.419.72"
$dCell.Style = $dStyle
$ws.Cells.Item(2, 5).Value = "  +0.02%  "

$dCell = $ws.Cells.Item(3, 4)
$dStyle = $dCell.Style
$dCell.Value = "'1.568.04"
$dCell.Style = $dStyle
$ws.Cells.Item(3, 5).Value = "  +0.19%  "

$ws.Cells.Item(4, 5).Value = "  +0.01%  "

$ws.Cells.Item(5, 5).Value = "  +0.03%  "

$ws.Cells.Item(6, 5).Value = "  +0.39%  "

$ws.Cells.Item(7, 5).Value = "  +1.74%  "

$dCell = $ws.Cells.Item(8, 4)
$dStyle = $dCell.Style
$dCell.Value = "'47.21"
$dCell.Style = $dStyle
$ws.Cells.Item(8, 5).Value = "  -2.81%  "

$dCell = $ws.Cells.Item(9, 4)
$dStyle = $dCell.Style
$dCell.Value = "'0.3295"
$dCell.Style = $dStyle
$ws.Cells.Item(9, 5).Value = "  -1.33%  "

$ws.Cells.Item(10, 5).Value = "  +2.44%  "

$dCell = $ws.Cells.Item(11, 4)
$dStyle = $dCell.Style
$dCell.Value = "'0.07439"
$dCell.Style = $dStyle
$ws.Cells.Item(11, 5).Value = "  +0.65%  "

$dCell = $ws.Cells.Item(12, 4)
$dStyle = $dCell.Style
$dCell.Value = "'1.001"
$dCell.Style = $dStyle
$ws.Cells.Item(12, 5).Value = "  +0.00%  "

$ws.Cells.Item(13, 5).Value = "  -0.31%  "

$dCell = $ws.Cells.Item(14, 4)
$dStyle = $dCell.Style
$dCell.Value = "'5.901"
$dCell.Style = $dStyle
$ws.Cells.Item(14, 5).Value = "  -0.64%  "

$ws.Cells.Item(15, 5).Value = "  -0.62%  "

$dCell = $ws.Cells.Item(16, 4)
$dStyle = $dCell.Style
$dCell.Value = "'1.569.54"
$dCell.Style = $dStyle
$ws.Cells.Item(16, 5).Value = "  +0.21%  "

$ws.Cells.Item(17, 5).Value = "  +0.49%  "

$dCell = $ws.Cells.Item(18, 4)
$dStyle = $dCell.Style
$dCell.Value = "'0.06707"
$dCell.Style = $dStyle
$ws.Cells.Item(18, 5).Value = "  +0.09%  "

$dCell = $ws.Cells.Item(19, 4)
$dStyle = $dCell.Style
$dCell.Value = "'86.65"
$dCell.Style = $dStyle
$ws.Cells.Item(19, 5).Value = "  -1.71%  "

$ws.Cells.Item(20, 5).Value = "  -0.04%  "

$dCell = $ws.Cells.Item(21, 4)
$dStyle = $dCell.Style
$dCell.Value = "'6.360"
$dCell.Style = $dStyle
$ws.Cells.Item(21, 5).Value = "  +0.30%  "

$dCell = $ws.Cells.Item(22, 4)
$dStyle = $dCell.Style
$dCell.Value = "'16.39"
$dCell.Style = $dStyle
$ws.Cells.Item(22, 5).Value = "  +1.51%  "

$ws.Cells.Item(23, 5).Value = "  -1.15%  "

$dCell = $ws.Cells.Item(24, 4)
$dStyle = $dCell.Style
$dCell.Value = "'22.402.15"
$dCell.Style = $dStyle
$ws.Cells.Item(24, 5).Value = "  -0.03%  "

$dCell = $ws.Cells.Item(25, 4)
$dStyle = $dCell.Style
$dCell.Value = "'2.353"
$dCell.Style = $dStyle
$ws.Cells.Item(25, 5).Value = "  -1.61%  "

$dCell = $ws.Cells.Item(26, 4)
$dStyle = $dCell.Style
$dCell.Value = "'2.597"
$dCell.Style = $dStyle
$ws.Cells.Item(26, 5).Value = "  +1.39%  "

$dCell = $ws.Cells.Item(27, 4)
$dStyle = $dCell.Style
$dCell.Value = "'151.05"
$dCell.Style = $dStyle
$ws.Cells.Item(27, 5).Value = "  +0.76%  "

$dCell = $ws.Cells.Item(28, 4)
$dStyle = $dCell.Style
$dCell.Value = "'19.49"
$dCell.Style = $dStyle
$ws.Cells.Item(28, 5).Value = "  +0.96%  "

$dCell = $ws.Cells.Item(29, 4)
$dStyle = $dCell.Style
$dCell.Value = "'4.933"
$dCell.Style = $dStyle
$ws.Cells.Item(29, 5).Value = "  -1.50%  "

$dCell = $ws.Cells.Item(30, 4)
$dStyle = $dCell.Style
$dCell.Value = "'124.16"
$dCell.Style = $dStyle
$ws.Cells.Item(30, 5).Value = "  +0.48%  "

$dCell = $ws.Cells.Item(31, 4)
$dStyle = $dCell.Style
$dCell.Value = "'1.745.36"
$dCell.Style = $dStyle
$ws.Cells.Item(31, 5).Value = "  +0.23%  "

$ws.Cells.Item(32, 5).Value = "  +1.66%  "

$dCell = $ws.Cells.Item(33, 4)
$dStyle = $dCell.Style
$dCell.Value = "'1.977"
$dCell.Style = $dStyle
$ws.Cells.Item(33, 5).Value = "  -1.32%  "

$dCell = $ws.Cells.Item(34, 4)
$dStyle = $dCell.Style
$dCell.Value = "'6.025"
$dCell.Style = $dStyle
$ws.Cells.Item(34, 5).Value = "  -1.15%  "

$dCell = $ws.Cells.Item(35, 4)
$dStyle = $dCell.Style
$dCell.Value = "'9.806"
$dCell.Style = $dStyle
$ws.Cells.Item(35, 5).Value = "  +0.31%  "

$dCell = $ws.Cells.Item(36, 4)
$dStyle = $dCell.Style
$dCell.Value = "'0.08278"
$dCell.Style = $dStyle
$ws.Cells.Item(36, 5).Value = "  +0.14%  "

$dCell = $ws.Cells.Item(37, 4)
$dStyle = $dCell.Style
$dCell.Value = "'0.02417"
$dCell.Style = $dStyle
$ws.Cells.Item(37, 5).Value = "  +0.54%  "

$dCell = $ws.Cells.Item(38, 4)
$dStyle = $dCell.Style
$dCell.Value = "'0.06356"
$dCell.Style = $dStyle
$ws.Cells.Item(38, 5).Value = "  -0.39%  "

$dCell = $ws.Cells.Item(39, 4)
$dStyle = $dCell.Style
$dCell.Value = "'1.287"
$dCell.Style = $dStyle
$ws.Cells.Item(39, 5).Value = "  -0.38%  "

$ws.Cells.Item(40, 5).Value = "  -1.25%  "

$dCell = $ws.Cells.Item(41, 4)
$dStyle = $dCell.Style
$dCell.Value = "'5.250"
$dCell.Style = $dStyle
$ws.Cells.Item(41, 5).Value = "  -1.36%  "

$dCell = $ws.Cells.Item(42, 4)
$dStyle = $dCell.Style
$dCell.Value = "'11.30"
$dCell.Style = $dStyle
$ws.Cells.Item(42, 5).Value = "  +1.71%  "

$dCell = $ws.Cells.Item(43, 4)
$dStyle = $dCell.Style
$dCell.Value = "'0.6136"
$dCell.Style = $dStyle
$ws.Cells.Item(43, 5).Value = "  +1.05%  "

$ws.Cells.Item(44, 5).Value = "  +0.97%  "

$dCell = $ws.Cells.Item(45, 4)
$dStyle = $dCell.Style
$dCell.Value = "'0.5976"
$dCell.Style = $dStyle
$ws.Cells.Item(45, 5).Value = "  +3.27%  "

$dCell = $ws.Cells.Item(46, 4)
$dStyle = $dCell.Style
$dCell.Value = "'3.752"
$dCell.Style = $dStyle
$ws.Cells.Item(46, 5).Value = "  -0.24%  "

$dCell = $ws.Cells.Item(47, 4)
$dStyle = $dCell.Style
$dCell.Value = "'2.023"
$dCell.Style = $dStyle
$ws.Cells.Item(47, 5).Value = "  +0.37%  "

$dCell = $ws.Cells.Item(48, 4)
$dStyle = $dCell.Style
$dCell.Value = "'124.62"
$dCell.Style = $dStyle
$ws.Cells.Item(48, 5).Value = "  +0.71%  "

$dCell = $ws.Cells.Item(49, 4)
$dStyle = $dCell.Style
$dCell.Value = "'1.192"
$dCell.Style = $dStyle
$ws.Cells.Item(49, 5).Value = "  -1.99%  "

$dCell = $ws.Cells.Item(50, 4)
$dStyle = $dCell.Style
$dCell.Value = "'0.07170"
$dCell.Style = $dStyle
$ws.Cells.Item(50, 5).Value = "  -0.38%  "

$dCell = $ws.Cells.Item(51, 4)
$dStyle = $dCell.Style
$dCell.Value = "'76.52"
$dCell.Style = $dStyle
$ws.Cells.Item(51, 5).Value = "  +1.25%  "
